# captest_prototype switch: add Bifaciality/StrucShd/BakMismatch inputs and
# sep/dayfirst/date_format CSV-reading config columns to the "PVsyst Runs" sheet,
# and nudge the saved scroll/selection state on "PVsyst Runs" and "Systems".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # PVsyst Runs
$ws2 = $wb.Worksheets.Item(2)   # Systems

# --- PVsyst Runs: insert a new column H ("Bifaciality") and push the old H:J -> I:K ---
$ws1.Columns.Item(8).Insert()

# --- Row 1 headers (string insertion order chosen to match the target shared-string table) ---
$ws1.Range("L1").Value = "sep"
$ws1.Range("M1").Value = "dayfirst"
$ws1.Range("L2").Value = ","
$ws1.Range("H1").Value = "Bifaciality"
$ws1.Range("N1").Value = "date_format"
$ws1.Range("N2").Value = "%m/%d/%y %H:%M"
$ws1.Range("O1").Value = "StrucShd"
$ws1.Range("P1").Value = "BakMismatch"

# --- Fill the rest of the new columns (H, L, M, N, O, P) for every data row ---
$bifaciality = @{2=0; 3=0.7; 4=0.7; 5=0.7; 6=0.7; 7=0; 8=0.7; 9=0.7}
$strucShd    = @{2=0; 3=0.05; 4=0.05; 5=0.05; 6=0.05; 7=0; 8=0.05; 9=0.05}
$bakMismatch = @{2=0; 3=0.1; 4=0.1; 5=0.1; 6=0.1; 7=0; 8=0.1; 9=0.1}

foreach ($r in 2..9) {
    $ws1.Range("H$r").Value = $bifaciality[$r]
    $ws1.Range("L$r").Value = ","
    $ws1.Range("M$r").Value = $false
    $ws1.Range("N$r").Value = "%m/%d/%y %H:%M"
    $ws1.Range("O$r").Value = $strucShd[$r]
    $ws1.Range("P$r").Value = $bakMismatch[$r]
}

# --- View/selection bookkeeping (matches the saved scroll + selection state) ---
# Set the non-active sheet's selection first so activating "PVsyst Runs" last leaves it the visible tab.
$ws2.Range("Y1:Z3").Select()
$ws1.Range("B1").Select()
$ws1.Range("O1").Select()
